$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.812.79'
$ws.Range("E2").Value = '  +0.27%  '
$ws.Range("D3").Value = '2.802.18'
$ws.Range("E3").Value = '  +0.80%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''354.52'
$ws.Range("E5").Value = '  -0.57%  '
$ws.Range("D6").Value = '''111.59'
$ws.Range("E6").Value = '  +2.23%  '
$ws.Range("D7").Value = '''0.558'
$ws.Range("E7").Value = '  +0.84%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").Value = '''0.636'
$ws.Range("E9").Value = '  +8.59%  '
$ws.Range("D10").Value = '''40.32'
$ws.Range("E10").Value = '  +1.59%  '
$ws.Range("E11").Value = '  -1.83%  '
$ws.Range("D12").Value = '''0.0840'
$ws.Range("E12").Value = '  -0.53%  '
$ws.Range("D13").Value = '''20.01'
$ws.Range("E13").Value = '  +2.57%  '
$ws.Range("D14").Value = '''7.78'
$ws.Range("E14").Value = '  +2.15%  '
$ws.Range("D15").Value = '3.243.32'
$ws.Range("E15").Value = '  +0.79%  '
$ws.Range("D16").Value = '2.802.03'
$ws.Range("E16").Value = '  +1.11%  '
$ws.Range("E17").Value = '  +1.59%  '
$ws.Range("D18").Value = '51.805.23'
$ws.Range("E18").Value = '  +0.28%  '
$ws.Range("D19").Value = '''7.67'
$ws.Range("E19").Value = '  +1.33%  '
$ws.Range("D20").Value = '''3.20'
$ws.Range("E20").Value = '  +3.84%  '
$ws.Range("D21").Value = '''13.71'
$ws.Range("E21").Value = '  +4.07%  '
$ws.Range("E22").Value = '  +1.11%  '
$ws.Range("D23").Value = '''70.50'
$ws.Range("E23").Value = '  +0.43%  '
$ws.Range("D24").Value = '''268.56'
$ws.Range("E24").Value = '  +0.44%  '
$ws.Range("E25").Value = '  +1.49%  '
$ws.Range("D26").Value = '''0.999'
$ws.Range("E26").Value = '  -0.13%  '
$ws.Range("D27").Value = '''26.18'
$ws.Range("E27").Value = '  -0.59%  '
$ws.Range("E28").Value = '  -2.39%  '
$ws.Range("D29").Value = '''38.90'
$ws.Range("E29").Value = '  +11.41%  '
$ws.Range("D30").Value = '''10.38'
$ws.Range("E31").Value = '  +3.18%  '
$ws.Range("D32").Value = '''6.15'
$ws.Range("E32").Value = '  -0.62%  '
$ws.Range("E33").Value = '  +0.74%  '
$ws.Range("D34").Value = '''5.63'
$ws.Range("E34").Value = '  +8.37%  '
$ws.Range("D35").Value = '''0.0881'
$ws.Range("E35").Value = '  +5.37%  '
$ws.Range("D36").Value = '''0.0444'
$ws.Range("E36").Value = '  -0.71%  '
$ws.Range("E37").Value = '  -0.04%  '
$ws.Range("E38").Value = '  +0.70%  '
$ws.Range("B39").Value = 'LidoDAOToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D39").Value = '''3.16'
$ws.Range("E39").Value = '  +0.90%  '
$ws.Range("B40").Value = 'ARBITRUM'
$ws.Range("C40").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D40").Value = '''2.01'
$ws.Range("E40").Value = '  +2.66%  '
$ws.Range("E41").Value = '  +1.26%  '
$ws.Range("D42").Value = '''2.50'
$ws.Range("E42").Value = '  -0.41%  '
$ws.Range("E43").Value = '  +1.31%  '
$ws.Range("D44").Value = '''120.68'
$ws.Range("E44").Value = '  +0.63%  '
$ws.Range("D45").Value = '''21.95'
$ws.Range("E45").Value = '  +1.47%  '
$ws.Range("E46").Value = '  +5.28%  '
$ws.Range("D47").Value = '2.113.97'
$ws.Range("E47").Value = '  +1.53%  '
$ws.Range("E48").Value = '  +6.28%  '
$ws.Range("E49").Value = '  +1.77%  '
$ws.Range("E50").Value = '  -1.06%  '
$ws.Range("E51").Value = '  +7.16%  '
